$p = $ppt.ActivePresentation

# Slide 3 ("Why is it Important?") - merge "Increased " + "Speeds" runs
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Characters(32, 17).Text = "Increased Speeds"

# Slide 4 ("Others' Work") - merge split quotation runs in each citation paragraph
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$tr4.Characters(1, 51).Text = [char]8220 + "The Art of Multiprocessor Programming" + [char]8221 + "- Herlihy & "
$tr4.Characters(66, 41).Text = [char]8220 + "Designing Concurrent Data Structures" + [char]8221 + " " + [char]8211 + " "
$tr4.Characters(128, 55).Text = [char]8220 + "Implementing Concurrent Data Objects" + [char]8221 + " " + [char]8211 + " Herlihy - 1993"

# Slide 7 ("Evaluation") - merge "Recorded iterations per " + "second against number of threads  "
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange
$tr7.Characters(92, 58).Text = "Recorded iterations per second against number of threads  "
